$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1045.2778
$ws.Range("I19").Value = 1043.1428
$ws.Range("J19").Value = 1046.6364
$ws.Range("K19").Value = 1043.1428
$ws.Range("L19").Value = 1046.6364
$ws.Range("M19").Value = -868.1428000000001
$ws.Range("N19").Value = -1396.6364

$ws.Range("H21").Value = 5085
$ws.Range("I21").Value = 170
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 170
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = 298
$ws.Range("N21").Value = -10936

$ws.Range("H23").Value = 5085
$ws.Range("I23").Value = 170
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 170
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = 64
$ws.Range("N23").Value = -10468

$ws.Range("H29").Value = 590.4
$ws.Range("J29").Value = 783.3333
$ws.Range("L29").Value = 2349.9999
$ws.Range("N29").Value = -2911.9999

$ws.Range("H38").Value = 2016226
$ws.Range("J38").Value = 250
$ws.Range("L38").Value = 750
$ws.Range("N38").Value = -1494

$ws.Range("H40").Value = 2041.6364
$ws.Range("I40").Value = 1957
$ws.Range("J40").Value = 2189.75
$ws.Range("K40").Value = 1957
$ws.Range("L40").Value = 2189.75
$ws.Range("M40").Value = -1782
$ws.Range("N40").Value = -2539.75

$ws.Range("H62").Value = 1970.2858
$ws.Range("I62").Value = 1965.5
$ws.Range("J62").Value = 1999
$ws.Range("K62").Value = 1965.5
$ws.Range("L62").Value = 1999
$ws.Range("M62").Value = -1341.5
$ws.Range("N62").Value = -3247

$ws.Range("H65").Value = 1970.2858
$ws.Range("I65").Value = 1965.5
$ws.Range("J65").Value = 1999
$ws.Range("K65").Value = 9827.5
$ws.Range("L65").Value = 9995
$ws.Range("M65").Value = -6707.5
$ws.Range("N65").Value = -16235

$ws.Range("H87").Value = 39153.6
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 39153.6
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 39153.6
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -41649.6

$ws.Range("H90").Value = 39153.6
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 39153.6
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 117460.8
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -129940.8

$ws.Range("H129").Value = 820.57776
$ws.Range("I129").Value = 557.6316
$ws.Range("J129").Value = 1012.7308
$ws.Range("K129").Value = 1672.8948
$ws.Range("L129").Value = 3038.1924
$ws.Range("M129").Value = 3327.1052
$ws.Range("N129").Value = -13038.1924

$ws.Range("H137").Value = 1391.6078
$ws.Range("I137").Value = 971.5
$ws.Range("K137").Value = 2914.5
$ws.Range("M137").Value = -364.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 504
$ws.Range("I19").Value = 504
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 504
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -275
$ws.Range("N19").ClearContents()

$ws.Range("H132").Value = 4524.577
$ws.Range("I132").Value = 4953.048
$ws.Range("J132").Value = 2725
$ws.Range("K132").Value = 14859.144
$ws.Range("L132").Value = 8175
$ws.Range("M132").Value = -12329.144
$ws.Range("N132").Value = -13235

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H120").Value = 37840
$ws.Range("J120").Value = 37840
$ws.Range("L120").Value = 37840
$ws.Range("N120").Value = -47516

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 41818
$ws.Range("J28").Value = 41818
$ws.Range("L28").Value = 41818
$ws.Range("N28").Value = -42308

$ws.Range("H132").Value = 2432.1428
$ws.Range("I132").Value = 1267.3334
$ws.Range("J132").Value = 3985.2222
$ws.Range("K132").Value = 3802.0002
$ws.Range("L132").Value = 11955.6666
$ws.Range("M132").Value = -1272.0002
$ws.Range("N132").Value = -17015.6666

$ws.Range("H134").Value = 1043.4
$ws.Range("I134").Value = 770.85297
$ws.Range("K134").Value = 2312.55891
$ws.Range("M134").Value = 222.4410899999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 17102.984
$ws.Range("J68").Value = 24618.814
$ws.Range("L68").Value = 73856.442
$ws.Range("N68").Value = -75478.442

$ws.Range("H69").Value = 2244.8572
$ws.Range("J69").Value = 2244.8572
$ws.Range("L69").Value = 6734.571599999999
$ws.Range("N69").Value = -8356.571599999999

$ws.Range("H71").Value = 17102.984
$ws.Range("J71").Value = 24618.814
$ws.Range("L71").Value = 221569.326
$ws.Range("N71").Value = -229681.326

$ws.Range("H72").Value = 2244.8572
$ws.Range("J72").Value = 2244.8572
$ws.Range("L72").Value = 20203.7148
$ws.Range("N72").Value = -28315.7148

$ws.Range("H131").Value = 853.8
$ws.Range("J131").Value = 853.8
$ws.Range("L131").Value = 2561.4
$ws.Range("N131").Value = -12641.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H46").Value = 8156.857
$ws.Range("I46").Value = 6200
$ws.Range("K46").Value = 6200
$ws.Range("M46").Value = -6044

$ws.Range("H57").Value = 19800
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 19800
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 19800
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -21440

$ws.Range("H122").Value = 2145.4211
$ws.Range("I122").Value = 2263.6667
$ws.Range("J122").Value = 1942.7142
$ws.Range("K122").Value = 6791.000100000001
$ws.Range("L122").Value = 5828.142599999999
$ws.Range("M122").Value = -4341.000100000001
$ws.Range("N122").Value = -10728.1426

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 6501.5
$ws.Range("J21").Value = 8333.333000000001
$ws.Range("L21").Value = 8333.333000000001
$ws.Range("N21").Value = -8681.333000000001

$ws.Range("H122").Value = 2900
$ws.Range("I122").Value = 2900
$ws.Range("K122").Value = 8700
$ws.Range("M122").Value = -6250

$ws.Range("H132").Value = 4319.25
$ws.Range("I132").Value = 4653.2383
$ws.Range("J132").Value = 3681.6365
$ws.Range("K132").Value = 13959.7149
$ws.Range("L132").Value = 11044.9095
$ws.Range("M132").Value = -11429.7149
$ws.Range("N132").Value = -16104.9095

$ws.Range("H136").Value = 1347.5883
$ws.Range("I136").Value = 990.4
$ws.Range("J136").Value = 1857.8572
$ws.Range("K136").Value = 2971.2
$ws.Range("L136").Value = 5573.571599999999
$ws.Range("M136").Value = -421.1999999999998
$ws.Range("N136").Value = -10673.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 4814.2856
$ws.Range("J20").Value = 4814.2856
$ws.Range("L20").Value = 4814.2856
$ws.Range("N20").Value = -5294.2856

$ws.Range("H132").Value = 2939.697
$ws.Range("I132").Value = 3558.8948
$ws.Range("J132").Value = 2099.3572
$ws.Range("K132").Value = 10676.6844
$ws.Range("L132").Value = 6298.071599999999
$ws.Range("M132").Value = -8146.6844
$ws.Range("N132").Value = -11358.0716

$ws.Range("H136").Value = 1378.1
$ws.Range("I136").Value = 687.1818
$ws.Range("J136").Value = 2222.5557
$ws.Range("K136").Value = 2061.5454
$ws.Range("L136").Value = 6667.6671
$ws.Range("M136").Value = 488.4546
$ws.Range("N136").Value = -11767.6671
Write-Output "Applied all cell updates."
